$wb = $excel.ActiveWorkbook

# --- Rename sheets: "Summary" -> "Claim", "Claims" -> "Submitted" ---
$wsClaim = $wb.Worksheets.Item("Summary")
$wsClaim.Name = "Claim"

$wsSubmitted = $wb.Worksheets.Item("Claims")
$wsSubmitted.Name = "Submitted"

# Fix up the Print_Titles defined name for the renamed "Submitted" sheet -
# Excel does not automatically rewrite the stored RefersTo formula text
# when the target sheet is renamed.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Submitted!Print_Titles") {
        $n.RefersTo = "=Submitted!`$1:`$1"
    }
}

# --- Update the title cell (locked, needs the sheet briefly unprotected) ---
$wsClaim.Unprotect()
$wsClaim.Range("A1").Value = "Claimbot"
$wsClaim.Protect()

# --- Update the demo username/password values to match their labels ---
$wsClaim.Range("B3").Value = "Username"
$wsClaim.Range("B4").Value = "Password"
